$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updated values (automatic electricity price update)
$ws.Range("A2").Value = 45885
$ws.Range("B2").Value = 118.82
$ws.Range("C2").Value = 114.32
$ws.Range("D2").Value = 109.69
$ws.Range("E2").Value = 105
$ws.Range("F2").Value = 104.28
$ws.Range("G2").Value = 105.35
$ws.Range("H2").Value = 108.43
$ws.Range("I2").Value = 111.19
$ws.Range("J2").Value = 103.67
$ws.Range("K2").Value = 60.93
$ws.Range("L2").Value = 26.28
$ws.Range("M2").Value = 3.8
$ws.Range("N2").Value = 2.01
$ws.Range("O2").Value = 2.01
$ws.Range("P2").Value = 0.66
$ws.Range("Q2").Value = 1.1
$ws.Range("R2").Value = 3.6
$ws.Range("S2").Value = 21.52
$ws.Range("T2").Value = 68.8
$ws.Range("U2").Value = 100
$ws.Range("V2").Value = 111.05
$ws.Range("W2").Value = 144
$ws.Range("X2").Value = 137
$ws.Range("Y2").Value = 117.65
$ws.Range("Z2").Value = 74.20999999999999

# AA2 (Slot_4h_max) stays "20h-24h"
$ws.Range("AB2").Value = 127.43

# AC2/AE2 slot labels swap, AD2/AF2 prices swap accordingly
$ws.Range("AC2").Value = "20h-22h"
$ws.Range("AD2").Value = 127.52
$ws.Range("AE2").Value = "22h-24h"
$ws.Range("AF2").Value = 127.32

# AG2 (Slot_min_price) stays "9h-18h"
